$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Tiva TM4C123GH6PM Microcontroller"
$ws.Range("B8").Value = "DATA SHEET"
$ws.Range("C8").Value = "Tiva TM4C123GH6PM"
$ws.Range("D8").Value = "tm4c123gh6pm2014.pdf"

$ws.Range("B8:C8").HorizontalAlignment = -4108

$ws.Range("D8").Select()
